$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 312 (high/low/close changed) ---
$ws.Cells.Item(312, 4).Value = 7.8715
$ws.Cells.Item(312, 5).Value = 7.8045
$ws.Cells.Item(312, 6).Value = 7.8712

# --- Append new rows 313-315, cloning row 312's formatting (A:G only) ---
$ws.Range("A312:G312").Copy($ws.Range("A313:G313"))
$ws.Range("A312:G312").Copy($ws.Range("A314:G314"))
$ws.Range("A312:G312").Copy($ws.Range("A315:G315"))

# Row 313
$ws.Cells.Item(313, 1).Value = 45170.33333333334
$ws.Cells.Item(313, 2).Value = "FX_IDC:USDGTQ"
$ws.Cells.Item(313, 3).Value = 7.8712
$ws.Cells.Item(313, 4).Value = 7.874
$ws.Cells.Item(313, 5).Value = 7.858
$ws.Cells.Item(313, 6).Value = 7.858
$ws.Cells.Item(313, 7).Value = 0

# Row 314
$ws.Cells.Item(314, 1).Value = 45201.375
$ws.Cells.Item(314, 2).Value = "FX_IDC:USDGTQ"
$ws.Cells.Item(314, 3).Value = 7.858
$ws.Cells.Item(314, 4).Value = 7.858
$ws.Cells.Item(314, 5).Value = 7.803
$ws.Cells.Item(314, 6).Value = 7.8325
$ws.Cells.Item(314, 7).Value = 0

# Row 315
$ws.Cells.Item(315, 1).Value = 45231.375
$ws.Cells.Item(315, 2).Value = "FX_IDC:USDGTQ"
$ws.Cells.Item(315, 3).Value = 7.8325
$ws.Cells.Item(315, 4).Value = 7.8325
$ws.Cells.Item(315, 5).Value = 7.824
$ws.Cells.Item(315, 6).Value = 7.8245
$ws.Cells.Item(315, 7).Value = 0
